# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# "Bad Drivers" table (rows 3-7)
# ----------------------------------------------------------------------

# Row 3 - label unchanged, only Critical Minutes / Good Roaming % change
$ws.Range("C3").Value = 721
$ws.Range("D3").Value = 94.8

# Row 4 - new driver version, new counts
$ws.Range("A4").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.18.2"
$ws.Range("B4").Value = 21
$ws.Range("C4").Value = 1871
$ws.Range("D4").Value = 98.8

# Row 5 - new driver version, new counts
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.60.1.2"
$ws.Range("B5").Value = 59
$ws.Range("C5").Value = 1693
$ws.Range("D5").Value = 98.90000000000001

# Row 6 - new driver version, new counts
$ws.Range("A6").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.30.1"
$ws.Range("B6").Value = 24
$ws.Range("C6").Value = 1136
$ws.Range("D6").Value = 98.90000000000001

# Row 7 - Totals
$ws.Range("B7").Value = 106
$ws.Range("C7").Value = 5421

# ----------------------------------------------------------------------
# "Good Drivers" table (rows 15-27)
# Pre-format the Driver Vintage column as text so date-looking strings
# are not auto-converted to Excel date serials.
# ----------------------------------------------------------------------
$ws.Range("E15:E27").NumberFormat = "@"

# Row 15
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B15").Value = 10661
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = "2022-08-29"

# Row 16
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B16").Value = 14239
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = "2022-05-23"

# Row 17
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B17").Value = 265400
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").Value = "2022-05-01"

# Row 18
$ws.Range("A18").Value = "Intel(R) Dual Band Wireless-AC 8265 - 22.30.0.11"
$ws.Range("B18").Value = 170510
$ws.Range("D18").Value = 99.90000000000001
$ws.Range("E18").Value = "2021-01-19"

# Row 19
$ws.Range("A19").Value = "Intel(R) Dual Band Wireless-AC 8265 - 22.0.1.1"
$ws.Range("B19").Value = 52096
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "2020-09-28"

# Row 20
$ws.Range("A20").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.11.3"
$ws.Range("B20").Value = 161874
$ws.Range("D20").Value = 100
$ws.Range("E20").Value = "2019-09-05"

# Row 21
$ws.Range("A21").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5"
$ws.Range("B21").Value = 143342
$ws.Range("D21").Value = 99.90000000000001
$ws.Range("E21").Value = "2019-08-25"

# Row 22
$ws.Range("A22").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.10.2"
$ws.Range("B22").Value = 20227
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = "2019-05-11"

# Row 23
$ws.Range("A23").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.9.1"
$ws.Range("B23").Value = 34065
$ws.Range("D23").Value = 100
$ws.Range("E23").Value = "2019-04-28"

# Row 24
$ws.Range("A24").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.8.1"
$ws.Range("B24").Value = 48540
$ws.Range("D24").Value = 100
$ws.Range("E24").Value = "2019-03-16"

# Row 25
$ws.Range("A25").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.5.2"
$ws.Range("B25").Value = 184564
$ws.Range("D25").Value = 99.90000000000001
$ws.Range("E25").Value = "2018-11-25"

# Row 26
$ws.Range("A26").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.50.0.4"
$ws.Range("B26").Value = 14221
$ws.Range("D26").Value = 100
$ws.Range("E26").Value = "2018-05-08"

# Row 27
$ws.Range("A27").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.30.1.2"
$ws.Range("B27").Value = 23765
$ws.Range("D27").Value = 100
$ws.Range("E27").Value = "2018-01-09"

# ----------------------------------------------------------------------
# Rows 28-30 previously held the three lowest-vintage "Good Driver" rows;
# the report now only has 13 good-driver rows (15-27), so clear them out.
# ----------------------------------------------------------------------
$ws.Range("A28:J30").Clear()
